# Método de Spline y HTML Sección 1
# Update the E (error) column values in the Newton's method iteration table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = "2.03049645777783"
    "D4"  = "1.1500233459578"
    "D5"  = "1.21236596596947"
    "D6"  = "0.885512699177286"
    "D7"  = "0.64559739903862"
    "D8"  = "0.0463501710601939"
    "D9"  = "0.0003239746094877"
    "D10" = "1.52778864409253e-08"
}

# Use a neighboring, already-text cell (D2) as the style template so the
# rewritten cells keep their original (default) style instead of picking
# up a new "quote prefix" style from the text-forcing apostrophe below.
$styleTemplate = $ws.Range("D2").Style

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces these numeric-looking strings to be stored
    # as text (matching the source inlineStr cells) instead of being
    # auto-coerced into numbers.
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = $styleTemplate
}
